$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.019.76"

$ws.Range("D3").Value = "'3.175.25"
$ws.Range("E3").Value = '  -0.84%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = "'603.66"
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("D6").Value = "'153.94"
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = "'3.175.78"
$ws.Range("E8").Value = '  -0.78%  '

$ws.Range("D9").Value = "'0.543"
$ws.Range("E9").Value = '  +1.62%  '

$ws.Range("E10").Value = '  -1.61%  '

$ws.Range("E11").Value = '  -7.70%  '

$ws.Range("D12").Value = "'0.506"
$ws.Range("E12").Value = '  -1.81%  '

$ws.Range("E13").Value = '  -2.81%  '

$ws.Range("D14").Value = "'38.20"
$ws.Range("E14").Value = '  -2.21%  '

$ws.Range("D15").Value = "'3.695.85"
$ws.Range("E15").Value = '  -0.81%  '

$ws.Range("D16").Value = "'66.058.66"
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").Value = "'7.34"
$ws.Range("E17").Value = '  -1.12%  '

$ws.Range("D18").Value = "'3.177.55"
$ws.Range("E18").Value = '  -0.73%  '

$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").Value = "'506.36"
$ws.Range("E20").Value = '  -1.11%  '

$ws.Range("D21").Value = "'15.24"
$ws.Range("E21").Value = '  -0.80%  '

$ws.Range("E22").Value = '  -2.09%  '

$ws.Range("D23").Value = "'7.99"
$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = "'14.74"
$ws.Range("E24").Value = '  -3.39%  '

$ws.Range("D25").Value = "'84.23"
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").Value = "'9.11"
$ws.Range("E28").Value = '  -3.18%  '

$ws.Range("E29").Value = '  +4.56%  '

$ws.Range("D30").Value = "'3.03"
$ws.Range("E30").Value = '  +5.26%  '

$ws.Range("D31").Value = "'6.99"
$ws.Range("E31").Value = '  +1.13%  '

$ws.Range("D32").Value = "'27.85"
$ws.Range("E32").Value = '  -1.68%  '

$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("E34").Value = '  -4.07%  '

$ws.Range("D35").Value = "'6.46"
$ws.Range("E35").Value = '  -1.70%  '

$ws.Range("D36").Value = "'509.44"
$ws.Range("E36").Value = '  +5.13%  '

$ws.Range("E37").Value = '  +0.37%  '

$ws.Range("D38").Value = "'0.0901"
$ws.Range("E38").Value = '  -0.42%  '

$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("D40").Value = "'0.0₃0709"
$ws.Range("E40").Value = '  +8.87%  '

$ws.Range("E41").Value = '  +4.05%  '

$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("E43").Value = '  -2.05%  '

$ws.Range("D44").Value = "'0.298"
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").Value = "'2.46"
$ws.Range("E45").Value = '  +1.44%  '

$ws.Range("D46").Value = "'2.826.61"
$ws.Range("E46").Value = '  -3.92%  '

$ws.Range("D47").Value = "'27.82"
$ws.Range("E47").Value = '  -2.24%  '

$ws.Range("E49").Value = '  +1.99%  '

$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("E51").Value = '  +2.31%  '
